$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 0.1468522516250963
$ws.Range("C7").Value = 0.6403690953631123
$ws.Range("D7").Value = 0.9600378145316305
$ws.Range("E7").Value = 0.9798151940706117
$ws.Range("F7").Value = 0.9824895536496444
$ws.Range("G7").Value = 36

$ws.Range("B8").Value = 0.1837530684691744
$ws.Range("C8").Value = 0.6684895077115939
$ws.Range("D8").Value = 0.9976342611431646
$ws.Range("E8").Value = 0.9988164301527906
$ws.Range("F8").Value = 0.9961014814234609
$ws.Range("G8").Value = 35

$ws.Range("B9").Value = 0.1757836139713026
$ws.Range("C9").Value = 0.7570773967417573
$ws.Range("D9").Value = 1.457116255293309
$ws.Range("E9").Value = 1.207110705483681
$ws.Range("F9").Value = 1.225267479434806
$ws.Range("G9").Value = 20

$ws.Range("B10").Value = -0.1910145113676734
$ws.Range("C10").Value = 0.6497709479478674
$ws.Range("D10").Value = 1.405754352160177
$ws.Range("E10").Value = 1.185645120666457
$ws.Range("F10").Value = 1.217938200125281
$ws.Range("G10").Value = 13

$ws.Range("B11").Value = 0.03995114952237531
$ws.Range("C11").Value = 0.4488742623708699
$ws.Range("D11").Value = 0.3692268520075374
$ws.Range("E11").Value = 0.6076403969516324
$ws.Range("F11").Value = 0.6778926515859446
$ws.Range("G11").Value = 5
